# Weekly fruit/vegetable price update: two new "Apio" price records were
# added to the daily logic subset. In the source data these new rows sit
# at the top of the date-ordered block (rows 419-420), pushing the
# existing records (419-440) down by two rows (to 421-442).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 419, shifting existing rows 419:440 down to 421:442
# (keeping all their original data/formatting intact).
$ws.Rows("419:420").Insert()

# New record 1 (row 419)
$ws.Range("A419").Value = 10
$ws.Range("B419").Value = "Vega Modelo de Temuco"
$ws.Range("C419").Value = "La Araucanía"
$ws.Range("D419").Value = 45008
$ws.Range("E419").Value = 9
$ws.Range("F419").Value = 100112017
$ws.Range("G419").Value = "Apio"
$ws.Range("H419").Value = "Americana (o)"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 140
$ws.Range("K419").Value = 12000
$ws.Range("L419").Value = 13000
$ws.Range("M419").Value = 12429
$ws.Range("N419").Value = "$/docena de matas"
$ws.Range("O419").Value = "Provincia del Elquí"
$ws.Range("P419").Value = 2072
$ws.Range("Q419").Value = 6
$ws.Range("R419").Value = "Hortaliza"

# New record 2 (row 420)
$ws.Range("A420").Value = 10
$ws.Range("B420").Value = "Vega Modelo de Temuco"
$ws.Range("C420").Value = "La Araucanía"
$ws.Range("D420").Value = 45008
$ws.Range("E420").Value = 9
$ws.Range("F420").Value = 100112017
$ws.Range("G420").Value = "Apio"
$ws.Range("H420").Value = "Americana (o)"
$ws.Range("I420").Value = "Segunda"
$ws.Range("J420").Value = 45
$ws.Range("K420").Value = 8000
$ws.Range("L420").Value = 8000
$ws.Range("M420").Value = 8000
$ws.Range("N420").Value = "$/docena de matas"
$ws.Range("O420").Value = "Provincia del Elquí"
$ws.Range("P420").Value = 1333
$ws.Range("Q420").Value = 6
$ws.Range("R420").Value = "Hortaliza"
